$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "7921854"
$ws.Range("D16").Value = "ROBERTO CARLOS ATENCIO CUELLO"
$ws.Range("E16").Value = "2008"
$ws.Range("F16").Value = 35112
$ws.Range("G16").Value = 807803

$ws.Range("C17").Value = "7921854"
$ws.Range("D17").Value = "ROBERTO CARLOS ATENCIO CUELLO"
$ws.Range("E17").Value = "2007"
$ws.Range("F17").Value = 35112
$ws.Range("G17").Value = 807803

$ws.Range("C18").Value = "1043312708"
$ws.Range("D18").Value = "JUAN LUIS CASTELLANO MARTINEZ"
$ws.Range("E18").Value = "1903"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = 828116

$ws.Range("C19").Value = "1043312708"
$ws.Range("D19").Value = "JUAN LUIS CASTELLANO MARTINEZ"
$ws.Range("E19").Value = "1902"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = 828116

$ws.Range("C20").Value = "1047370395"
$ws.Range("D20").Value = "OCTAVIO ENRIQUE RAUDALES GARRIDO"
$ws.Range("E20").Value = "2008"
$ws.Range("F20").Value = 35112
$ws.Range("G20").Value = 807803

$ws.Range("C21").Value = "1047370395"
$ws.Range("D21").Value = "OCTAVIO ENRIQUE RAUDALES GARRIDO"
$ws.Range("E21").Value = "2007"
$ws.Range("F21").Value = 35112
$ws.Range("G21").Value = 807803
